# Update "RunMode1" metadata sheet with the next batch of Cocci results
# (20220919-Cocci-10427Updt), replacing the previous Salm-13035 batch
# values in columns A (Result ID), E (Lab Sample ID) and T (Cartridge ID)
# for data rows 2-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labSampleId = "20220919-Cocci-10427Updt"
$cartridgeId = "TestCartridge0427"

$resultIds = @(
    "A1892501",
    "A1892502",
    "A1892503",
    "A1892504",
    "A1892505",
    "A1892506",
    "A1892507",
    "A1892508",
    "A1892509",
    "A1892510",
    "A1892511",
    "A1892512"
)

for ($i = 0; $i -lt $resultIds.Length; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $resultIds[$i]
    $ws.Range("E$row").Value = $labSampleId
    $ws.Range("T$row").Value = $cartridgeId
}
